$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 3.1
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 3.6
$ws.Range("AJ5").Value = 9.5
$ws.Range("BB5").Value = 151
